$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '72.249.02'
$ws.Range("E2").Value = '  +0.46%  '

$ws.Range("D3").Value = '2.706.70'
$ws.Range("E3").Value = '  +2.85%  '

$ws.Range("E4").Value = '  +0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '600.41'
$ws.Range("E5").Value = '  -1.03%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '175.68'
$ws.Range("E6").Value = '  -2.15%  '

$ws.Range("E7").Value = '  +0.04%  '

$ws.Range("E8").Value = '  -0.21%  '

$ws.Range("D9").Value = '2.706.57'
$ws.Range("E9").Value = '  +2.85%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.169'
$ws.Range("E10").Value = '  +0.17%  '

$ws.Range("E11").Value = '  +2.62%  '

$ws.Range("E12").Value = '  +1.96%  '

$ws.Range("E13").Value = '  -0.17%  '

$ws.Range("D14").Value = '3.203.55'
$ws.Range("E14").Value = '  +1.55%  '

$ws.Range("E15").Value = '  -0.26%  '

$ws.Range("D16").Value = '72.083.91'
$ws.Range("E16").Value = '  +0.44%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.35'
$ws.Range("E17").Value = '  -0.70%  '

$ws.Range("D18").Value = '2.711.68'
$ws.Range("E18").Value = '  +3.77%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.28'
$ws.Range("E19").Value = '  +6.83%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.13'
$ws.Range("E20").Value = '  +2.15%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '373.09'
$ws.Range("E21").Value = '  -3.16%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.19'
$ws.Range("E22").Value = '  +0.64%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.02'
$ws.Range("E23").Value = '  +2.01%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '72.44'
$ws.Range("E24").Value = '  -0.50%  '

$ws.Range("E25").Value = '  -0.07%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '4.38'
$ws.Range("E26").Value = '  -1.82%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.82'
$ws.Range("E27").Value = '  -0.82%  '

$ws.Range("D28").Value = '2.847.37'
$ws.Range("E28").Value = '  +2.93%  '

$ws.Range("E29").Value = '  +0.13%  '

$ws.Range("D30").Value = '0.0₃0992'
$ws.Range("E30").Value = '  +2.93%  '

$ws.Range("E31").Value = '  +0.77%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '509.39'
$ws.Range("E32").Value = '  -6.74%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.31'
$ws.Range("E33").Value = '  -1.65%  '

$ws.Range("E34").Value = '  -0.11%  '

$ws.Range("E35").Value = '  +0.03%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '164.09'
$ws.Range("E36").Value = '  -1.15%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '19.72'
$ws.Range("E37").Value = '  +2.50%  '

$ws.Range("E38").Value = '  -0.22%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.39'
$ws.Range("E39").Value = '  -0.22%  '

$ws.Range("E40").Value = '  -4.71%  '

$ws.Range("E41").Value = '  -3.24%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.07'
$ws.Range("E42").Value = '  +0.65%  '

$ws.Range("E43").Value = '  +0.02%  '

$ws.Range("E44").Value = '  +0.89%  '

$ws.Range("E45").Value = '  -2.38%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '157.05'
$ws.Range("E46").Value = '  +3.89%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '39.51'
$ws.Range("E47").Value = '  +0.77%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.566'
$ws.Range("E48").Value = '  +5.77%  '

$ws.Range("E49").Value = '  +2.96%  '

$ws.Range("E50").Value = '  +5.08%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0764'
$ws.Range("E51").Value = '  +0.98%  '
